$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update join-date column (D) from date-serial values to plain year numbers,
# and clear the date number formatting on those cells.
$ws.Range("D2:D4").Style = "Normal"
$ws.Range("D2").Value = 2005
$ws.Range("D3").Value = 2006
$ws.Range("D4").Value = 2007

# Move the active selection to G10 (was L4).
$ws.Range("G10").Select()
